# Generate Report for Handoff
# Updates the localization-status workbook with a freshly generated handoff
# package: new GUID-named markdown/xlf files, refreshed handoff timestamps,
# and cleared handback (target) info since the new handoff has not yet been
# translated/returned.

$wb = $excel.ActiveWorkbook

$oldGuid = "24659008-6d08-43ac-8072-5cfee2563ff9"
$newGuid = "b041d819-2578-4f69-b3d6-99b18f6e30eb"
$oldHash = "013fa165bbd107f6d1205b71446064b63bcc4385"
$newHash = "053063ce0ff96a4f8d62efcd49f997ba3ed21cf7"

$newHandoffDate = "2016-08-20 01:02:46"
$newZhHandoffDatetime = "2016-08-20 01:02:42"
$emptyDatetime = "0001-01-01 00:00:00"

$sourceMdName = $newGuid + ".md"
$zhXlfName = $newGuid + "." + $newHash + ".zh-cn.xlf"
$deXlfName = $newGuid + "." + $newHash + ".de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $sourceMdName
$wsOverview.Range("G2").Value = $newHandoffDate

# The hyperlink target URL is unchanged (still points at the original commit
# in the ol-test0 repo) - only the displayed text is refreshed to the new
# handoff file name.
$overviewLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34c5794040b23001f712814c6baf6c9feae82c1c/e2e/" + $oldGuid + ".md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewLinkAddress, "", "", "e2e\" + $sourceMdName)

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $sourceMdName
$wsZh.Range("G2").Value = $zhXlfName
$wsZh.Range("H2").Value = $newZhHandoffDatetime
$wsZh.Range("I2").Value = ""
$wsZh.Range("I2").Style = "Normal"
$wsZh.Range("J2").Value = ""
$wsZh.Range("K2").Value = $emptyDatetime

$zhLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34c5794040b23001f712814c6baf6c9feae82c1c/e2e/" + $oldGuid + ".md"
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhLinkAddress, "", "", $sourceMdName)

$wsZh.Columns.Item(9).ColumnWidth = 17.75
$wsZh.Columns.Item(10).ColumnWidth = 20.75

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $sourceMdName
$wsDe.Range("G2").Value = $deXlfName
$wsDe.Range("H2").Value = $newHandoffDate
$wsDe.Range("I2").Value = ""
$wsDe.Range("I2").Style = "Normal"
$wsDe.Range("J2").Value = ""
$wsDe.Range("K2").Value = $emptyDatetime

$deLinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34c5794040b23001f712814c6baf6c9feae82c1c/e2e/" + $oldGuid + ".md"
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deLinkAddress, "", "", $sourceMdName)

$wsDe.Columns.Item(9).ColumnWidth = 17.75
$wsDe.Columns.Item(10).ColumnWidth = 20.75
